$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 18: RR vs KKR (row 30) - fill in player scores
$ws.Range("E30").Value = 20
$ws.Range("H30").Value = 30
$ws.Range("K30").Value = 60
$ws.Range("N30").Value = 80
$ws.Range("Q30").Value = 40
$ws.Range("T30").Value = 70
$ws.Range("W30").Value = 50
$ws.Range("Z30").Value = 100
$ws.Range("AC30").Value = 0

$wb.Application.CalculateFullRebuild()
